# Update NATMI LR-pair output (Wnt1-Fzd4) with new TPM-based results.
# The "ECs" sending-cluster rows are removed entirely, leaving only the
# "FAPs" sending-cluster rows (previously rows 5-7, now rows 2-4), whose
# derived-specificity / weight columns are refreshed with values computed
# from the new TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 2-4 (Sending cluster = "ECs"); this shifts the former
# rows 5-7 (Sending cluster = "FAPs") up to become the new rows 2-4.
$ws.Range("A2:A4").EntireRow.Delete()

# Row 2 (FAPs -> ECs)
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Wnt1"
$ws.Cells.Item(2, 3).Value = "Fzd4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.3532066666666667
$ws.Cells.Item(2, 8).Value = 1.05962
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 28.22405966666667
$ws.Cells.Item(2, 14).Value = 84.672179
$ws.Cells.Item(2, 15).Value = 0.3816548478108986
$ws.Cells.Item(2, 16).Value = 0.3816548478108986
$ws.Cells.Item(2, 17).Value = 9.968926034664443
$ws.Cells.Item(2, 18).Value = 89.72033431198
$ws.Cells.Item(2, 19).Value = 0.3816548478108986
$ws.Cells.Item(2, 20).Value = 0.3816548478108986

# Row 3 (FAPs -> FAPs)
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Wnt1"
$ws.Cells.Item(3, 3).Value = "Fzd4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3532066666666667
$ws.Cells.Item(3, 8).Value = 1.05962
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 19.768727
$ws.Cells.Item(3, 14).Value = 59.306181
$ws.Cells.Item(3, 15).Value = 0.2673191094302723
$ws.Cells.Item(3, 16).Value = 0.2673191094302723
$ws.Cells.Item(3, 17).Value = 6.982446167913333
$ws.Cells.Item(3, 18).Value = 62.84201551122
$ws.Cells.Item(3, 19).Value = 0.2673191094302723
$ws.Cells.Item(3, 20).Value = 0.2673191094302723

# Row 4 (FAPs -> MuSCs)
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Wnt1"
$ws.Cells.Item(4, 3).Value = "Fzd4"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.3532066666666667
$ws.Cells.Item(4, 8).Value = 1.05962
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 25.95900466666667
$ws.Cells.Item(4, 14).Value = 77.877014
$ws.Cells.Item(4, 15).Value = 0.351026042758829
$ws.Cells.Item(4, 16).Value = 0.351026042758829
$ws.Cells.Item(4, 17).Value = 9.168893508297778
$ws.Cells.Item(4, 18).Value = 82.52004157468001
$ws.Cells.Item(4, 19).Value = 0.351026042758829
$ws.Cells.Item(4, 20).Value = 0.351026042758829
